# Re-applies the refreshed "cryptos" price/volume snapshot (GitHub Actions bot).
# Numeric-looking Price values are written with a leading apostrophe so Excel
# stores them as literal text (matching the source data, e.g. "61.364.49" /
# "550.02") instead of silently coercing them to floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range('D2').Value = '61.364.49'
$ws.Range('E2').Value = '  +0.20%  '

# row 3
$ws.Range('D3').Value = '2.380.91'
$ws.Range('E3').Value = '  +0.08%  '

# row 4
$ws.Range('E4').Value = '  +0.04%  '

# row 5
$ws.Range('D5').Value = '''550.02'
$ws.Range('E5').Value = '  +0.29%  '

# row 6
$ws.Range('D6').Value = '''139.28'
$ws.Range('E6').Value = '  -1.61%  '

# row 7
$ws.Range('E7').Value = '  -0.04%  '

# row 8
$ws.Range('E8').Value = '  -0.86%  '

# row 9
$ws.Range('D9').Value = '2.382.24'
$ws.Range('E9').Value = '  +0.19%  '

# row 10
$ws.Range('E10').Value = '  +2.11%  '

# row 11
$ws.Range('E11').Value = '  +1.37%  '

# row 12
$ws.Range('E12').Value = '  +0.71%  '

# row 13
$ws.Range('E13').Value = '  +0.74%  '

# row 14
$ws.Range('D14').Value = '''25.24'
$ws.Range('E14').Value = '  -0.23%  '

# row 15
$ws.Range('E15').Value = '  +1.27%  '

# row 16
$ws.Range('D16').Value = '61.277.65'
$ws.Range('E16').Value = '  +0.15%  '

# row 17
$ws.Range('D17').Value = '2.373.76'
$ws.Range('E17').Value = '  -0.19%  '

# row 18
$ws.Range('E18').Value = '  +2.40%  '

# row 19
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '''321.63'
$ws.Range('E19').Value = '  +0.97%  '

# row 20
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '''4.15'
$ws.Range('E20').Value = '  +0.78%  '

# row 21
$ws.Range('D21').Value = '''6.74'
$ws.Range('E21').Value = '  +1.04%  '

# row 22
$ws.Range('E22').Value = '  +0.05%  '

# row 23
$ws.Range('D23').Value = '''64.29'
$ws.Range('E23').Value = '  +0.75%  '

# row 24
$ws.Range('D24').Value = '''1.71'
$ws.Range('E24').Value = '  -9.59%  '

# row 25
$ws.Range('D25').Value = '''8.62'
$ws.Range('E25').Value = '  +5.50%  '

# row 26
$ws.Range('D26').Value = '''8.18'
$ws.Range('E26').Value = '  +1.36%  '

# row 27
$ws.Range('D27').Value = '''511.68'
$ws.Range('E27').Value = '  -3.04%  '

# row 28
$ws.Range('D28').Value = '0.0₃0895'
$ws.Range('E28').Value = '  -3.45%  '

# row 29
$ws.Range('E29').Value = '  +3.28%  '

# row 30
$ws.Range('E30').Value = '  -2.88%  '

# row 31
$ws.Range('E31').Value = '  -0.11%  '

# row 32
$ws.Range('E32').Value = '  -1.90%  '

# row 33
$ws.Range('E33').Value = '  -0.02%  '

# row 34
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '''4.70'
$ws.Range('E34').Value = '  +0.99%  '

# row 35
$ws.Range('B35').Value = 'Stacks'
$ws.Range('C35').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D35').Value = '''1.91'
$ws.Range('E35').Value = '  +3.20%  '

# row 36
$ws.Range('E36').Value = '  +0.70%  '

# row 37
$ws.Range('E37').Value = '  +1.13%  '

# row 38
$ws.Range('E38').Value = '  +2.32%  '

# row 39
$ws.Range('E39').Value = '  +4.53%  '

# row 40
$ws.Range('E40').Value = '  -0.11%  '

# row 41
$ws.Range('D41').Value = '''41.21'
$ws.Range('E41').Value = '  +1.87%  '

# row 42
$ws.Range('D42').Value = '''151.17'
$ws.Range('E42').Value = '  +7.57%  '

# row 43
$ws.Range('D43').Value = '''2.16'
$ws.Range('E43').Value = '  +0.68%  '

# row 44
$ws.Range('D44').Value = '''3.59'
$ws.Range('E44').Value = '  -0.22%  '

# row 45
$ws.Range('E45').Value = '  +1.09%  '

# row 46
$ws.Range('D46').Value = '''19.47'
$ws.Range('E46').Value = '  -2.91%  '

# row 47
$ws.Range('E47').Value = '  +0.79%  '

# row 48
$ws.Range('E48').Value = '  +0.32%  '

# row 49
$ws.Range('E49').Value = '  -0.40%  '

# row 50
$ws.Range('E50').Value = '  +0.47%  '

# row 51
$ws.Range('D51').Value = '''16.81'

Write-Output "Updated cryptos list"
